$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the data rows (2-9), formulas + values, keep them blank ---
$ws.Range("A2:D9").ClearContents()

# --- Copy the plain (non-bold, theme-colored, no-fill) header format from B1
#     onto A1 so both headers start from the same base font before we
#     diverge them with Bold/Color below (A1 currently carries the old
#     "Inconsolata + yellow fill" formatting we want to get rid of). ---
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Turn the header row into plain text values (drop the IMPORTRANGE
#     formulas), pulled from shared strings. ---
$ws.Range("A1").Value = "Заказчик"
$ws.Range("B1").Value = "Юр.лицо"
$ws.Range("C1").Value = "Город"
$ws.Range("D1").Value = "Адрес"

# --- Bold the whole header row. ---
$ws.Range("A1:D1").Font.Bold = $true

# --- A1 gets an explicit black (non-theme) colour + re-normalised name so
#     it ends up on its own distinct font (vs. the bold/theme font shared
#     by B1:D1). ---
$ws.Range("A1").Font.Color = 0
$ws.Range("A1").Font.Name = "Calibri"

# --- A1 gets the white fill back (it already existed in the palette). ---
$ws.Range("A1").Interior.Color = 16777215

# --- Thin box border around every header cell. ---
$ws.Range("A1:D1").Borders.LineStyle = 1

# --- Header row no longer needs the taller (16.2pt) formula-driven height. ---
$ws.Rows(1).AutoFit()

# --- Column sizing: widen City / add Address column. ---
$ws.Columns("C").ColumnWidth = 18.666666666666668
$ws.Columns("D").ColumnWidth = 22.333333333333332

# --- Selection + print orientation. ---
$ws.Range("A2").Select()
$ws.PageSetup.Orientation = 1
